# "cho thue xe, tien tru dan vao tai khoan"
# Add a new "money" column (header + starting balance) to the user sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header in L1 (shared string "money") and seed value 100 in L2.
$ws.Range("L1").Value = "money"
$ws.Range("L2").Value = 100

# Leave the selection where Excel left it after adding the column.
$ws.Range("L3").Select()
